$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old standalone "Overall average..." row (was row 24); it gets
# re-created at row 25 once the new "vph" rows are inserted above it.
$ws.Rows("24").Delete()

# New row 23: per-direction throughput converted from cars/minute to vph
# (cars/hour), i.e. the row-22 averages multiplied by 60.
$ws.Range("A23").Value = "vph"
$ws.Range("B23").Formula = "=B22*60"
$ws.Range("C23:M23").Formula = "=C22*60"
$ws.Range("B23:M23").NumberFormat = "0"

# New row 24: directional-flow labels (e.g. "north_south") matching the
# column headers in row 1, one per movement column B..M.
$ws.Range("B24").Value = "north_south"
$ws.Range("C24").Value = "north_east"
$ws.Range("D24").Value = "north_west"
$ws.Range("E24").Value = "south_north"
$ws.Range("F24").Value = "south_west"
$ws.Range("G24").Value = "south_east"
$ws.Range("H24").Value = "west_east"
$ws.Range("I24").Value = "west_north"
$ws.Range("J24").Value = "west_south"
$ws.Range("K24").Value = "east_west"
$ws.Range("L24").Value = "east_south"
$ws.Range("M24").Value = "east_north"
$ws.Range("B24:M24").NumberFormat = "0"

# Row 25: the overall-average summary line, now bold to match the sheet's
# other summary cells.
$ws.Range("A25").Value = "Overall average throughput per minute = 35.45 cars per minute"
$ws.Range("A25").Font.Bold = $true

[void]$ws.Range("D25").Select()
